$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$reqFraco = "LOT2028 -  Tecnologia de Processos Fermentativos  (Requisito fraco)`n"
$indConjunto = "LOT2052 -  Tecnologia de Bebidas Experimental  (Indicação de Conjunto)`n"

# Swap the two requisito lines: row 23 becomes the "Indicação de Conjunto" entry
# and row 24 becomes the "Requisito fraco" entry (both columns B and C mirror each other).
$ws.Range("B23").Value = $indConjunto
$ws.Range("C23").Value = $indConjunto
$ws.Range("B24").Value = $reqFraco
$ws.Range("C24").Value = $reqFraco
